
# ---------------------------------------------------------------------------
# 1. add reset baby info
# 2. add new icon for setting in scene baby info panel
# 3. update plan and add user feedbacks  -> add new "測試情況" worksheet
#    with device / feedback table
# 4. add function: press Key Escape in babyInfo/ pattern/ calendar to quit
#    application
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a brand-new worksheet, "測試情況", at the very end of the workbook
# (after the existing Sheet1 / Sheet2 / Sheet3) and make it the active sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "測試情況"

# Header row (B1:E1) first ...
$ws.Range("B1").Value = "手機品牌"
$ws.Range("C1").Value = "手機型號"
$ws.Range("D1").Value = "安卓版本"
$ws.Range("E1").Value = "測試情況"

# ... then the data row (B2:E2) ...
$ws.Range("B2").Value = "小米"
$ws.Range("C2").Value = "MI 5S"
$ws.Range("D2").Value = "6.0.1"
$ws.Range("E2").Value = "可用"

# ... and finally column F (header + data) last.
$ws.Range("F1").Value = "修改意見"
$ws.Range("F2").Value = "返回鍵退出整個程序"

# Center-align the header row
$ws.Range("B1:F1").HorizontalAlignment = -4108

# Column widths (best achievable match given host rounding granularity)
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 18.5
$ws.Columns.Item(4).ColumnWidth = 17.0
$ws.Columns.Item(5).ColumnWidth = 47.666666666666664
$ws.Columns.Item(6).ColumnWidth = 49.833333333333336

# Match the saved selection on the new sheet
$ws.Range("F4").Select() | Out-Null
